$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# publication_id column
$ws.Range("A6").Value = 121
$ws.Range("A7").Value = 1343

# title column
$ws.Range("B6").Value = "Fatal arrhythmias associated with genetic variants in type 2 ryanodine receptor channel gene"
$ws.Range("B7").Value = "Novel mutations in arrhythmogenic right ventricular cardiomyopathy from Indian population"

# first_author column
$ws.Range("C6").Value = "Horie, M."
$ws.Range("C7").Value = "Pamuru, PR"

# doi column
$ws.Range("D6").Value = "10.1007/s10840-018-0338-y"
$ws.Range("D7").Value = "10.4103/0971-6866.86182"

# Apply the red fill used throughout the table to the newly added publication id cells
$ws.Range("A6:A7").Interior.Color = 255

# Leave the selection where the user ended up after the edit
$ws.Range("M15:M16").Select() | Out-Null
